# Updated cryptos list on Thu Apr 25 10:58:50 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 47 / 48 swap: ThetaToken moves up to row 47, USDe moves down to row 48 ---
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -5.50%  "

$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.02%  "

# --- Price (D) and Volume(1h) (E) updates for all remaining rows ---
$ws.Range("D2").Value = "63.606.45"
$ws.Range("E2").Value = "  -4.17%  "

$ws.Range("D3").Value = "3.089.75"
$ws.Range("E3").Value = "  -5.35%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "606.86"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").Value = "144.16"
$ws.Range("E6").Value = "  -8.56%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "3.086.04"
$ws.Range("E8").Value = "  -5.50%  "

$ws.Range("E9").Value = "  -5.23%  "

$ws.Range("E10").Value = "  -8.24%  "

$ws.Range("E11").Value = "  -10.85%  "

$ws.Range("D12").Value = "0.465"
$ws.Range("E12").Value = "  -6.45%  "

$ws.Range("E13").Value = "  -9.21%  "

$ws.Range("D14").Value = "34.88"
$ws.Range("E14").Value = "  -10.66%  "

$ws.Range("D15").Value = "3.602.76"
$ws.Range("E15").Value = "  -5.06%  "

$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("D17").Value = "63.668.32"
$ws.Range("E17").Value = "  -4.17%  "

$ws.Range("D18").Value = "3.092.71"
$ws.Range("E18").Value = "  -5.13%  "

$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  -9.35%  "

$ws.Range("D20").Value = "472.37"
$ws.Range("E20").Value = "  -6.61%  "

$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  -6.08%  "

$ws.Range("E22").Value = "  -7.79%  "

$ws.Range("E23").Value = "  -6.12%  "

$ws.Range("D24").Value = "13.42"
$ws.Range("E24").Value = "  -8.66%  "

$ws.Range("D25").Value = "82.79"
$ws.Range("E25").Value = "  -4.97%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("E27").Value = "  -9.62%  "

$ws.Range("E28").Value = "  -10.27%  "

$ws.Range("E29").Value = "  -11.46%  "

$ws.Range("D30").Value = "6.63"
$ws.Range("E30").Value = "  -5.60%  "

$ws.Range("E31").Value = "  -13.63%  "

$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("E33").Value = "  -5.94%  "

$ws.Range("D34").Value = "25.92"
$ws.Range("E34").Value = "  -7.50%  "

$ws.Range("E35").Value = "  -4.82%  "

$ws.Range("D36").Value = "5.87"
$ws.Range("E36").Value = "  -9.42%  "

$ws.Range("D37").Value = "52.04"
$ws.Range("E37").Value = "  -6.62%  "

$ws.Range("D38").Value = "0.0₃0721"
$ws.Range("E38").Value = "  -8.63%  "

$ws.Range("D39").Value = "452.08"
$ws.Range("E39").Value = "  -8.95%  "

$ws.Range("D40").Value = "2.87"
$ws.Range("E40").Value = "  -15.03%  "

$ws.Range("D41").Value = "0.0389"
$ws.Range("E41").Value = "  -7.96%  "

$ws.Range("E42").Value = "  -8.31%  "

$ws.Range("D43").Value = "8.27"
$ws.Range("E43").Value = "  -6.48%  "

$ws.Range("D44").Value = "2.808.99"
$ws.Range("E44").Value = "  -6.53%  "

$ws.Range("E45").Value = "  -10.71%  "

$ws.Range("D46").Value = "2.21"
$ws.Range("E46").Value = "  -12.69%  "

$ws.Range("D49").Value = "25.78"
$ws.Range("E49").Value = "  -11.03%  "

$ws.Range("E50").Value = "  -5.95%  "

$ws.Range("D51").Value = "117.40"
$ws.Range("E51").Value = "  -2.48%  "
